# Update "想去人数" (column F) counts across the four worksheets to the
# newly scraped values, as produced by the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> list of (row, new value) pairs for column F.
$updates = @{
    "展览" = @{
        3  = 171
        6  = 1181
        7  = 2274
        8  = 2193
        10 = 624
        12 = 1722
        13 = 422
        16 = 315
        17 = 261
        18 = 1639
        19 = 292
        20 = 1349
        21 = 764
        22 = 292
        24 = 12474
        25 = 12531
        27 = 722
        29 = 274
        31 = 427
        32 = 1947
        36 = 637
    }
    "演出" = @{
        4 = 44
        8 = 121
    }
    "本地生活" = @{
        2 = 79
        3 = 126
    }
    "全部类型" = @{
        3  = 79
        4  = 171
        7  = 1181
        8  = 2274
        9  = 2193
        11 = 624
        12 = 126
        14 = 1722
        15 = 422
        20 = 315
        21 = 44
        22 = 261
        23 = 1639
        24 = 292
        25 = 1349
        26 = 764
        27 = 292
        30 = 12474
        31 = 12531
        33 = 722
        35 = 274
        37 = 427
        40 = 1947
        42 = 121
        46 = 637
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $newValue = $rows[$row]
        $ws.Cells.Item([int]$row, 6).Value = $newValue
    }
}
